# Update the "week" column (B2:B54) by adding 26 to each existing value,
# then set the active selection to L13 (no multi-cell selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 54; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value2 = $cell.Value2 + 26
}

$ws.Range("L13").Select()
